# Apply the "disconnected_elements" diagnostic layout:
#   B1 = 0   (bold, thin boxed border, centered/top aligned)
#   A2 = 0   (same style as B1)
#   B2 = "disconnected_elements" (plain, shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full style once on B1: bold font + thin border all around + center/top alignment
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160

# A2 gets the same value + style; copy the finished format from B1 so no
# intermediate style is ever recorded for A2.
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# B2 is a plain label cell (goes into the shared strings table)
$ws.Range("B2").Value = "disconnected_elements"
